# Edit script: applies the changes described by the diff to po-cond01.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the cont1_same / cont1_opp rows (original rows 2 and 3).
#    Everything below shifts up automatically, carrying formatting with it.
$ws.Rows("2:3").Delete()

# 2. Insert a new column for "logStep" right after "rampLin" (original col K),
#    i.e. before the old "postStimBlankT" column (original col L).
$ws.Columns("L").Insert()
$ws.Range("L1").Value = "logStep"
$ws.Range("L2:L10").Value = 1

# 3. Insert a new column for "maskRamp" right after "maskRR" (now col N after
#    the previous insert), i.e. before the old "odtLoc" column.
$ws.Columns("O").Insert()
$ws.Range("O1").Value = "maskRamp"
$ws.Range("O2:O10").Value = 0

# 4. stimT: 1000 -> 2000 for all data rows (col F)
$ws.Range("F2:F10").Value = 2000

# 5. odtTilt: 1.5 -> 2 for all data rows (now col Q after the column inserts)
$ws.Range("Q2:Q10").Value = 2

# 6. nRevs: 12 -> 20 for all data rows (now col W after the column inserts)
$ws.Range("W2:W10").Value = 20

# 7. Update the view state to match (scrolled right, selection on V11)
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("V11").Select()
